{"js": "// Update the 25 two-digit multiplication equations in the body table.\n// Each old equation string is unique in the document, so we can safely\n// find-and-replace each one independently via a case-sensitive search.\nconst replacements = [\n  [\"64\u00d789=5696\", \"51\u00d736=1836\"],\n  [\"13\u00d753=689\", \"26\u00d780=2080\"],\n  [\"31\u00d734=1054\", \"58\u00d743=2494\"],\n  [\"59\u00d736=2124\", \"72\u00d758=4176\"],\n  [\"25\u00d731=775\", \"51\u00d738=1938\"],\n  [\"99\u00d716=1584\", \"75\u00d785=6375\"],\n  [\"22\u00d777=1694\", \"84\u00d751=4284\"],\n  [\"44\u00d749=2156\", \"43\u00d725=1075\"],\n  [\"92\u00d766=6072\", \"78\u00d740=3120\"],\n  [\"67\u00d761=4087\", \"12\u00d790=1080\"],\n  [\"34\u00d741=1394\", \"70\u00d784=5880\"],\n  [\"17\u00d723=391\", \"99\u00d733=3267\"],\n  [\"60\u00d729=1740\", \"67\u00d783=5561\"],\n  [\"21\u00d777=1617\", \"19\u00d764=1216\"],\n  [\"70\u00d783=5810\", \"26\u00d747=1222\"],\n  [\"62\u00d763=3906\", \"68\u00d785=5780\"],\n  [\"82\u00d755=4510\", \"93\u00d760=5580\"],\n  [\"33\u00d763=2079\", \"96\u00d762=5952\"],\n  [\"36\u00d728=1008\", \"47\u00d714=658\"],\n  [\"27\u00d739=1053\", \"96\u00d786=8256\"],\n  [\"73\u00d715=1095\", \"39\u00d736=1404\"],\n  [\"75\u00d733=2475\", \"40\u00d781=3240\"],\n  [\"50\u00d733=1650\", \"87\u00d776=6612\"],\n  [\"52\u00d790=4680\", \"88\u00d773=6424\"],\n  [\"15\u00d755=825\", \"85\u00d742=3570\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Could not find text to replace: ${oldText}`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the 25 two-digit multiplication equations in the body table.\n# Each old equation string is unique in the document, so a plain\n# Find/Replace (whole text, case-sensitive) for each pair is safe.\n\n$wdReplaceAll   = 2\n$wdFindContinue = 1\n\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"64\u00d789=5696\", \"51\u00d736=1836\"),\n    @(\"13\u00d753=689\",  \"26\u00d780=2080\"),\n    @(\"31\u00d734=1054\", \"58\u00d743=2494\"),\n    @(\"59\u00d736=2124\", \"72\u00d758=4176\"),\n    @(\"25\u00d731=775\",  \"51\u00d738=1938\"),\n    @(\"99\u00d716=1584\", \"75\u00d785=6375\"),\n    @(\"22\u00d777=1694\", \"84\u00d751=4284\"),\n    @(\"44\u00d749=2156\", \"43\u00d725=1075\"),\n    @(\"92\u00d766=6072\", \"78\u00d740=3120\"),\n    @(\"67\u00d761=4087\", \"12\u00d790=1080\"),\n    @(\"34\u00d741=1394\", \"70\u00d784=5880\"),\n    @(\"17\u00d723=391\",  \"99\u00d733=3267\"),\n    @(\"60\u00d729=1740\", \"67\u00d783=5561\"),\n    @(\"21\u00d777=1617\", \"19\u00d764=1216\"),\n    @(\"70\u00d783=5810\", \"26\u00d747=1222\"),\n    @(\"62\u00d763=3906\", \"68\u00d785=5780\"),\n    @(\"82\u00d755=4510\", \"93\u00d760=5580\"),\n    @(\"33\u00d763=2079\", \"96\u00d762=5952\"),\n    @(\"36\u00d728=1008\", \"47\u00d714=658\"),\n    @(\"27\u00d739=1053\", \"96\u00d786=8256\"),\n    @(\"73\u00d715=1095\", \"39\u00d736=1404\"),\n    @(\"75\u00d733=2475\", \"40\u00d781=3240\"),\n    @(\"50\u00d733=1650\", \"87\u00d776=6612\"),\n    @(\"52\u00d790=4680\", \"88\u00d773=6424\"),\n    @(\"15\u00d755=825\",  \"85\u00d742=3570\")\n)\n\nforeach ($pair in $pairs) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($pair[0], $false, $true, $false, $false, $false, $true, $wdFindContinue, $false, $pair[1], $wdReplaceAll)\n}\n"}
